# Fill in the two blank tracking rows (Index 7 and 8 -> sheet rows 9 and 10)
# on the "Spieltabelle" sheet with the new session data, and move the
# active-cell selection to M17 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (Index 7) ---------------------------------------------------
$ws.Range("B9").Value = "Cashgame"
$ws.Range("C9").Value = "sc.ch"
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 0.8
$ws.Range("F9").Value = 0.01
$ws.Range("G9").Value = 45949
$ws.Range("H9").Value = 0.72569444444444442
$ws.Range("I9").Value = 45949
$ws.Range("J9").Value = 0.76736111111111116
$ws.Range("K9").Value = 1.5
$ws.Range("L9").Value = 3.11
$ws.Range("M9").Value = 3.11
$ws.Range("N9").Value = 3
$ws.Range("O9").Value = 0.31
$ws.Range("P9").Value = "Nichts"
$ws.Range("Q9").Value = "Durchschnittliche Karten mit (ein Full House)"
$ws.Range("R9").Value = "Sehr viel limping, C-Bets haben meisten funktioniert"

# --- Row 10 (Index 8) ---------------------------------------------------
$ws.Range("B10").Value = "Cashgame"
$ws.Range("C10").Value = "sc.ch"
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 0.8
$ws.Range("F10").Value = 0.01
$ws.Range("G10").Value = 45949
$ws.Range("H10").Value = 0.76736111111111116
$ws.Range("I10").Value = 45949
$ws.Range("J10").Value = 0.86111111111111116
$ws.Range("K10").Value = 1.5
$ws.Range("L10").Value = 2.38
$ws.Range("M10").Value = 2.5
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 0.4
$ws.Range("P10").Value = "Nichts"
$ws.Range("Q10").Value = "3mal 2Paare und gegen Brodway verloren"
$ws.Range("R10").Value = "Sehr viel limping, Viele 3way Pots, sehr viel Tilting (shanesas)"

# --- Final selection, matching the saved cursor position ---------------
$ws.Range("M17").Select()
